$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "30.131.59"
Set-TextValue $ws.Range("E2") "  -0.17%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.873.90"
Set-TextValue $ws.Range("E3") "  -1.02%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.9992"
Set-TextValue $ws.Range("E4") "  -0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "242.81"
Set-TextValue $ws.Range("E5") "  -1.70%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.9993"
Set-TextValue $ws.Range("E6") "  +0.07%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4876"
Set-TextValue $ws.Range("E7") "  -2.23%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.2890"
Set-TextValue $ws.Range("E8") "  -1.51%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.06574"
Set-TextValue $ws.Range("E9") "  -1.20%  "

# Row 10
Set-TextValue $ws.Range("D10") "1.867.29"
Set-TextValue $ws.Range("E10") "  -1.21%  "

# Row 11
Set-TextValue $ws.Range("E11") "  -4.14%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.07192"
Set-TextValue $ws.Range("E12") "  -0.18%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.6635"
Set-TextValue $ws.Range("E13") "  -1.75%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "4.931"
Set-TextValue $ws.Range("E14") "  +1.68%  "

# Row 15
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D15") "85.87"
Set-TextValue $ws.Range("E15") "  -0.28%  "

# Row 16
Set-TextValue $ws.Range("D16") "30.069.18"
Set-TextValue $ws.Range("E16") "  -0.39%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D17") "0.000007774"
Set-TextValue $ws.Range("E17") "  -2.62%  "

# Row 18
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D18") "0.9995"
Set-TextValue $ws.Range("E18") "  +0.15%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.73"
Set-TextValue $ws.Range("E19") "  -0.91%  "

# Row 20
$ws.Range("B20").Value = "BinanceUSD"
$ws.Range("C20").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D20") "1.007"
Set-TextValue $ws.Range("E20") "  +0.77%  "

# Row 21
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D21") "2.113.08"
Set-TextValue $ws.Range("E21") "  -1.07%  "

# Row 22
Set-TextValue $ws.Range("D22") "4.747"
Set-TextValue $ws.Range("E22") "  -0.69%  "

# Row 23
Set-TextValue $ws.Range("D23") "5.873"
Set-TextValue $ws.Range("E23") "  +4.28%  "

# Row 24
Set-TextValue $ws.Range("D24") "9.152"
Set-TextValue $ws.Range("E24") "  -0.16%  "

# Row 25
Set-TextValue $ws.Range("D25") "151.73"
Set-TextValue $ws.Range("E25") "  +2.29%  "

# Row 26
Set-TextValue $ws.Range("D26") "143.39"
Set-TextValue $ws.Range("E26") "  +6.07%  "

# Row 27
Set-TextValue $ws.Range("D27") "16.90"
Set-TextValue $ws.Range("E27") "  +0.59%  "

# Row 28
Set-TextValue $ws.Range("D28") "1.877"
Set-TextValue $ws.Range("E28") "  -3.33%  "

# Row 29
Set-TextValue $ws.Range("E29") "  +0.87%  "

# Row 30
Set-TextValue $ws.Range("D30") "4.190"
Set-TextValue $ws.Range("E30") "  -0.40%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.08777"
Set-TextValue $ws.Range("E31") "  +0.16%  "

# Row 32
Set-TextValue $ws.Range("D32") "3.986"
Set-TextValue $ws.Range("E32") "  +0.56%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.05144"
Set-TextValue $ws.Range("E33") "  -0.16%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.7153"
Set-TextValue $ws.Range("E34") "  +1.03%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.106"
Set-TextValue $ws.Range("E35") "  -1.71%  "

# Row 36
Set-TextValue $ws.Range("D36") "2.665"
Set-TextValue $ws.Range("E36") "  +0.00%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.01843"
Set-TextValue $ws.Range("E37") "  +11.04%  "

# Row 38
Set-TextValue $ws.Range("D38") "2.670"
Set-TextValue $ws.Range("E38") "  -3.97%  "

# Row 39
Set-TextValue $ws.Range("D39") "2.157"
Set-TextValue $ws.Range("E39") "  -3.65%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.9262"
Set-TextValue $ws.Range("E40") "  -1.71%  "

# Row 41
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D41") "103.85"
Set-TextValue $ws.Range("E41") "  +0.74%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D42") "0.4233"
Set-TextValue $ws.Range("E42") "  +0.51%  "

# Row 43
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D43") "0.9985"
Set-TextValue $ws.Range("E43") "  +0.18%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D44") "5.750"
Set-TextValue $ws.Range("E44") "  -5.52%  "

# Row 45
Set-TextValue $ws.Range("D45") "7.432"
Set-TextValue $ws.Range("E45") "  -1.26%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.1280"
Set-TextValue $ws.Range("E46") "  +1.23%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.05714"
Set-TextValue $ws.Range("E47") "  -0.23%  "

# Row 48
Set-TextValue $ws.Range("D48") "32.76"
Set-TextValue $ws.Range("E48") "  -0.12%  "

# Row 49
Set-TextValue $ws.Range("D49") "8.251"
Set-TextValue $ws.Range("E49") "  -0.41%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.3753"
Set-TextValue $ws.Range("E50") "  +0.41%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.342"
Set-TextValue $ws.Range("E51") "  -0.13%  "
